$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (style) from row 876 into new rows 877-891
$srcFmt = $ws.Range("A876:C876")
for ($r = 877; $r -le 891; $r++) {
    $dstFmt = $ws.Range("A" + $r + ":C" + $r)
    $srcFmt.Copy($dstFmt) | Out-Null
}

# Set language column (A) for all new rows
for ($r = 877; $r -le 891; $r++) {
    $ws.Cells.Item($r, 1).Value = "cs"
}

# Set label/translation values in the exact original authoring order
# so shared-string table indices line up with the source workbook
$ws.Cells.Item(877, 2).Value = "lab.vape.leaks.0"
$ws.Cells.Item(877, 3).Value = "Žádné"
$ws.Cells.Item(878, 2).Value = "lab.vape.leaks.1"
$ws.Cells.Item(878, 3).Value = "Drobné"
$ws.Cells.Item(879, 3).Value = "Větší"
$ws.Cells.Item(881, 2).Value = "lab.vape.leaks.2"
$ws.Cells.Item(879, 2).Value = "lab.vape.leaks.3"
$ws.Cells.Item(880, 3).Value = "Totální"
$ws.Cells.Item(880, 2).Value = "lab.vape.leaks.4"
$ws.Cells.Item(881, 3).Value = "Občasné"
$ws.Cells.Item(882, 2).Value = "lab.vape.dryhit.0"
$ws.Cells.Item(882, 3).Value = "Žádné"
$ws.Cells.Item(883, 2).Value = "lab.vape.dryhit.1"
$ws.Cells.Item(883, 3).Value = "Občas"
$ws.Cells.Item(884, 2).Value = "lab.vape.dryhit.2"
$ws.Cells.Item(884, 3).Value = "Velmi často"
$ws.Cells.Item(885, 2).Value = "lab.vape.airflow.0"
$ws.Cells.Item(885, 3).Value = "Utažený"
$ws.Cells.Item(886, 2).Value = "lab.vape.airflow.1"
$ws.Cells.Item(886, 3).Value = "Mírně utažený"
$ws.Cells.Item(887, 2).Value = "lab.vape.airflow.2"
$ws.Cells.Item(887, 3).Value = "Volný"
$ws.Cells.Item(888, 2).Value = "lab.vape.airflow.3"
$ws.Cells.Item(888, 3).Value = "Zcela otevřený"
$ws.Cells.Item(889, 2).Value = "lab.vape.juice.0"
$ws.Cells.Item(889, 3).Value = "Utažený"
$ws.Cells.Item(890, 2).Value = "lab.vape.juice.1"
$ws.Cells.Item(890, 3).Value = "Mírně otevřený"
$ws.Cells.Item(891, 2).Value = "lab.vape.juice.2"
$ws.Cells.Item(891, 3).Value = "Zcela otevřený"

$ws.Activate() | Out-Null
$ws.Range("B889").Select() | Out-Null
